# JsonForBuyerSeller.xlsx edit script
# "Added procedure notes to excel"
#
# 1. Rename Sheet1 -> Practice
# 2. Add a new "Notes" worksheet after Practice, populate it with Q&A notes
# 3. Update the price/diff table on Practice: several numeric "diff" cells
#    become the literal string "Deduced"
# 4. Add a "Query" note cell (L8) and an empty formatted cell (L17) to Practice
# 5. Rename buyer/seller/diff -> buyInCountry/sellInCountry/profit in the
#    CONCAT() formulas (L20:N23)
# 6. Widen a few columns, update selection/active sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the first sheet
# ---------------------------------------------------------------------
$practice = $wb.Worksheets.Item(1)
$practice.Name = "Practice"

# ---------------------------------------------------------------------
# 2. Add the Notes sheet right after Practice
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Add($null, $practice)
$notes.Name = "Notes"

# ---------------------------------------------------------------------
# 3. Practice: rewrite the "diff" table -> several cells become "Deduced"
# ---------------------------------------------------------------------
$practice.Range("D4").Value = "Deduced"
$practice.Range("D5").Value = "Deduced"
$practice.Range("E5").Value = "Deduced"
$practice.Range("D6").Value = "Deduced"
$practice.Range("E6").Value = "Deduced"
$practice.Range("F6").Value = "Deduced"

# ---------------------------------------------------------------------
# 4. Practice: new note cells
# ---------------------------------------------------------------------
$practice.Range("L8").Value = "Query: If I have to find what to buy in USA: I will query for buyer or seller as USA and …. "
$practice.Range("L8").WrapText = $true
$practice.Rows(8).RowHeight = 28.8

$practice.Range("L17").WrapText = $true

# ---------------------------------------------------------------------
# 5. Practice: update the CONCAT formulas (buyer/seller/diff -> buyInCountry/sellInCountry/profit)
# ---------------------------------------------------------------------
$practice.Range("L20").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(C3)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(D3)-COLUMN($C$2),1,1),", profit:",D3,",product:", $C$2,"}")'
$practice.Range("M20").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(D3)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(E3)-COLUMN($C$2),1,1),", profit:",E3,",product:", $C$2,"}")'
$practice.Range("N20").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(E3)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(F3)-COLUMN($C$2),1,1),", profit:",F3,",product:", $C$2,"}")'

$practice.Range("L21").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(C4)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(D4)-COLUMN($C$2),1,1),", profit:",D4,",product:", $C$2,"}")'
$practice.Range("M21").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(D4)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(E4)-COLUMN($C$2),1,1),", profit:",E4,",product:", $C$2,"}")'
$practice.Range("N21").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(E4)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(F4)-COLUMN($C$2),1,1),", profit:",F4,",product:", $C$2,"}")'

$practice.Range("L22").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(C5)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(D5)-COLUMN($C$2),1,1),", profit:",D5,",product:", $C$2,"}")'
$practice.Range("M22").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(D5)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(E5)-COLUMN($C$2),1,1),", profit:",E5,",product:", $C$2,"}")'
$practice.Range("N22").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(E5)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(F5)-COLUMN($C$2),1,1),", profit:",F5,",product:", $C$2,"}")'

$practice.Range("L23").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(C6)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(D6)-COLUMN($C$2),1,1),", profit:",D6,",product:", $C$2,"}")'
$practice.Range("M23").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(D6)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(E6)-COLUMN($C$2),1,1),", profit:",E6,",product:", $C$2,"}")'
$practice.Range("N23").Formula = '=_xlfn.CONCAT("{buyInCountry:",OFFSET($C$2,ROW(E6)-ROW($C$2),0,1,1),",sellInCountry:",OFFSET($C$2,0,COLUMN(F6)-COLUMN($C$2),1,1),", profit:",F6,",product:", $C$2,"}")'

# ---------------------------------------------------------------------
# 6. Practice: widen columns that now hold longer text
# ---------------------------------------------------------------------
$practice.Columns("C:C").ColumnWidth = 35.8
$practice.Columns("L:L").ColumnWidth = 65.6
$practice.Columns("M:M").ColumnWidth = 62.8
$practice.Columns("N:N").ColumnWidth = 61.92

# ---------------------------------------------------------------------
# Populate the Notes sheet
# ---------------------------------------------------------------------
$notes.Range("A1").Value = "What concrete questions do we want to answer, how to answer and how well does our data model fare in answering those"
$notes.Range("A1").WrapText = $true
$notes.Rows(1).RowHeight = 28.8

$notes.Range("A6").Value = "Question"
$notes.Range("B6").Value = "Answer procedure"

$notes.Range("A7").Value = "Given a product (say Samsung s9) which country has the lowest and highest price for it "
$notes.Range("A7").VerticalAlignment = -4108
$notes.Range("B7").Value = "PreCondition: No repeating buyerCountry and sellerCountry, individually repeatable but not together`nGiven my data model looks like this {BuyInCountry: , SellInCountry: , profit:, product:} I will query by product=%Desired_product% and profit > 0 and order by highest to lowest, I will make another query with product=%Desired_product% and profit < 0 ordered by lowest to highest, switch the sign on profit, switch the buyer and seller and merge the lists."
$notes.Range("B7").VerticalAlignment = -4160
$notes.Range("B7").WrapText = $true
$notes.Rows(7).RowHeight = 86.4

$notes.Range("A8").Value = "Given I am travelling to somewhere (say USA) from somewhere (say India) what are the items I should buy while going to there and coming back from there in order to generate maximum profit"
$notes.Range("A8").HorizontalAlignment = -4131
$notes.Range("A8").VerticalAlignment = -4108
$notes.Range("A8").WrapText = $true
$notes.Range("B8").Value = "I will first visit the site when I am in India. I will select destination as USA and select 'seller' as my trader type, the system first makes a call with buyer as India and seller as US where profit > 0 and arrange the list by descending order of profit, it then makes a call as Seller as india and Buyer as US and profit < 0 and it arranges in decreasing order (highest negative to lowest negative), it then flips the sign and the buyer and seller (so buyer becomes India and seller becomes US)`nI will then visit the site from US where the process reverses"
$notes.Range("B8").WrapText = $true
$notes.Rows(8).RowHeight = 86.4

$notes.Range("A9").Value = "Given a country (say India) I want to find out what item I should buy here that will generate maximum profit for me and where should I sell it."
$notes.Range("A9").WrapText = $true
$notes.Range("B9").Value = "I will query on buyer as India or seller as India and for buyer as India call I will order by descending order in profit and for seller as India call I filter on profit < 0, I will order by ascending order of profit (highest negative to lowest negative) and flip the sign and the buyer and seller"
$notes.Range("B9").WrapText = $true
$notes.Rows(9).RowHeight = 43.2

$notes.Range("A12").Value = "Spacial requirements"

$notes.Range("A13").Value = "each entry comprising of {BuyInCountry: , SellInCountry: , profit:, product:} has 100 chars (400 Bytes) so 10 items for 4KB, 100 items for 40KB, 1000 items for 4MB, 10000 items for 40MB, 1 Lakh for 4 GB"
$notes.Range("A13").WrapText = $true
$notes.Rows(13).RowHeight = 43.2

$notes.Range("A14").Value = "25 Products * 20 Countries * 20 Countries - So given 4 GB of space I can get top 25 products for 20 countries"
$notes.Range("A14").WrapText = $true
$notes.Rows(14).RowHeight = 28.8

$notes.Columns("A:A").ColumnWidth = 80.88671875
$notes.Columns("B:B").ColumnWidth = 82.77734375

# ---------------------------------------------------------------------
# Selection / active sheet
# ---------------------------------------------------------------------
$practice.Range("L2").Select()
$notes.Activate()
$notes.Range("A13").Select()
